# Drugscreen template: add a dedicated "conc_condition" concentration column
# (so max_conc no longer needs to be redefined elsewhere) and rename the
# placeholder "Null" header to "condition_string".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column in front of the old "Concentration" column (J) and
# give it a header + the same default value (0.01) as its neighbours.
$ws.Range("J1").EntireColumn.Insert()

# Rename the "Null" placeholder string used in I2 to "condition_string".
$ws.Range("I2").Value = "condition_string"

# Populate the freshly inserted column J.
$ws.Range("J1").Value = "conc_condition"
$ws.Range("J2").Value = 0.01

# Move the active selection to K2 (the old "Concentration" data cell).
$ws.Range("K2").Select()
